# Fix template sheet name:
# the first worksheet was named "Отчет по классификаторам" (leftover from a
# classifier-report template) but this workbook is the load-test report
# template, so rename it to match. Excel sheet names are capped at 31
# characters, so "Отчет по нагрузочному тестированию" is truncated to
# "Отчет по нагрузочному тестирова".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Отчет по классификаторам")
$ws.Name = "Отчет по нагрузочному тестирова"
